$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CJ20")

# Insert a new row at row 33 (the 260 A frame size was missing from the
# table, between the existing 225 A and 265 A rows). This shifts rows
# 33-44 down to 34-45.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row with the new data point, matching the
# formatting/style already used by the rest of the table (s="2", the
# centered style applied via column A's style and used throughout
# A3:E44).
$ws.Range("A33:E33").HorizontalAlignment = -4108
$ws.Cells.Item(33, 1).Value = "CJ20"
$ws.Cells.Item(33, 2).Value = 0.66
$ws.Cells.Item(33, 3).Value = "3P"
$ws.Cells.Item(33, 4).Value = 260
$ws.Cells.Item(33, 5).Value = "DIN-Rail"

# Update the selection to reflect where the user left off after the
# edit: the new row's amperage cell.
$ws.Range("D34").Select()
